# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the two substantive changes from the target commit:
#
#   1. The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") is
#      switched from the deck's custom table style to a built-in
#      PowerPoint table style ({47030691-432D-4747-A2AD-122929E1CF28}).
#
#   2. The presentation's applied theme ("Integral" / Red Violet) is
#      swapped for the standard Office theme's colour scheme
#      ("Office Theme" / Office), which is the same set of colours
#      that lived in the deck's other (otherwise-unused) theme part.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Table style on slide 5's table (Shape 2: "Google Shape;122;p17")
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{47030691-432D-4747-A2AD-122929E1CF28}")
}

# ---------------------------------------------------------------------
# 2. Theme colours: swap the "Integral" (Red Violet) scheme currently
#    applied to the deck for the standard "Office Theme" scheme.
#    The 12 DrawingML theme colours, in order, are:
#      dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# ---------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    # PowerPoint COM RGB colors are stored BGR-packed (like the RGB() macro)
    $colorScheme.Item($i).RGB = ($b * 0x10000) + ($g * 0x100) + $r
}
